# Adds result-modal language strings for: PC Verify, Device Acquisition,
# and Volatile Data Gather (plus a "Points/Penalty" key and a
# "Digital Investigation Suite" title key).
#
# The "en" sheet is a simple Key/Value table (col A = key, col B = value).
# Before this edit the last row (111) holds percent / "Percent:".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- PC Verify result modal -------------------------------------------------
$ws.Range("A112").Value = "pc_verify_result"
$ws.Range("B112").Value = "PC Verify Result"

$ws.Range("A113").Value = "pc_verify_check_network"
$ws.Range("B113").Value = "Network Cable Unplugged"

$ws.Range("A114").Value = "pc_verify_check_power"
$ws.Range("B114").Value = "Check PC Power Status"

$ws.Range("A115").Value = "pc_verify_capture_screen"
$ws.Range("B115").Value = "Monitor Active Screen Pictured"

# --- Volatile data gather result modal --------------------------------------
$ws.Range("A116").Value = "volatile_gather_result"
$ws.Range("B116").Value = "Volatile Data Acquisition Result"

$ws.Range("A117").Value = "volatile_gather_result_order"
$ws.Range("A118").Value = "volatile_gather_result_player_order"
$ws.Range("B117").Value = "Recommended Order"
$ws.Range("B118").Value = "Your Order"

# --- Device acquisition result modal ----------------------------------------
$ws.Range("A119").Value = "device_gather_result"
$ws.Range("B119").Value = "Device Acquisition Result"

# --- Misc: points/penalty key, inserted just above the old last row --------
$ws.Rows.Item(111).Insert()
$ws.Range("B111").Value = "Points/Penalty"
$ws.Range("A111").Value = "points_penalty"

# --- Misc: suite title key, appended as the new last row --------------------
$ws.Range("B121").Value = "Digital Investigation Suite"
$ws.Range("A121").Value = "digital_investigation_suite"

$ws.Range("A121").Select() | Out-Null
